$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Refresh timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 13:35"

# Alemania (row 10) - updated case counts
$ws.Range("B10").Value = 172626
$ws.Range("C10").Value = 50
$ws.Range("E10").Value = 17765

# Corea del Sur / Sudafrica swap places in the ranking (row 43 <-> row 44)
$ws.Range("A43").Value = "Sudafrica"
$ws.Range("B43").Value = 11350
$ws.Range("C43").Value = 698
$ws.Range("D43").Value = 4357
$ws.Range("E43").Value = 6787
$ws.Range("F43").Value = 77
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 206
$ws.Range("A44").Value = "Corea del Sur"
$ws.Range("B44").Value = 10936
$ws.Range("C44").Value = 27
$ws.Range("D44").Value = 9670
$ws.Range("E44").Value = 1008
$ws.Range("F44").Value = 55
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 258

# Australia (row 53) - updated active/recovered counts
$ws.Range("D53").Value = 6229
$ws.Range("E53").Value = 644

# Kazajistan / Barein swap places in the ranking (row 59 <-> row 60)
$ws.Range("A59").Value = "Barein"
$ws.Range("B59").Value = 5409
$ws.Range("C59").Value = 173
$ws.Range("D59").Value = 2182
$ws.Range("E59").Value = 3218
$ws.Range("F59").Value = 5
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 9
$ws.Range("A60").Value = "Kazajistan"
$ws.Range("B60").Value = 5279
$ws.Range("C60").Value = 72
$ws.Range("D60").Value = 2108
$ws.Range("E60").Value = 3139
$ws.Range("F60").Value = 33
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 32

# Sri Lanka / Libano swap places in the ranking (row 105 <-> row 106)
$ws.Range("A105").Value = "Libano"
$ws.Range("B105").Value = 870
$ws.Range("C105").Value = 11
$ws.Range("D105").Value = 234
$ws.Range("E105").Value = 610
$ws.Range("F105").Value = 3
$ws.Range("H105").Value = 26
$ws.Range("A106").Value = "Sri Lanka"
$ws.Range("B106").Value = 869
$ws.Range("C106").Value = 6
$ws.Range("D106").Value = 366
$ws.Range("E106").Value = 494
$ws.Range("F106").Value = 1
$ws.Range("H106").Value = 9

# Benin moves above Montenegro / Republica del Chad (rows 134-136 rotate)
$ws.Range("A134").Value = "Benin"
$ws.Range("B134").Value = 327
$ws.Range("C134").Value = 8
$ws.Range("D134").Value = 76
$ws.Range("E134").Value = 249
$ws.Range("F134").Value = 0
$ws.Range("H134").Value = 2
$ws.Range("A135").Value = "Montenegro"
$ws.Range("B135").Value = 324
$ws.Range("D135").Value = 294
$ws.Range("E135").Value = 21
$ws.Range("F135").Value = 2
$ws.Range("H135").Value = 9
$ws.Range("A136").Value = "Republica del Chad"
$ws.Range("B136").Value = 322
$ws.Range("D136").Value = 53
$ws.Range("E136").Value = 238
$ws.Range("H136").Value = 31

# Camboya / Uganda swap places in the ranking (row 158 <-> row 159)
$ws.Range("A158").Value = "Uganda"
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 55
$ws.Range("E158").Value = 67
$ws.Range("F158").Value = 0
$ws.Range("A159").Value = "Camboya"
$ws.Range("B159").Value = 122
$ws.Range("D159").Value = 121
$ws.Range("E159").Value = 1
$ws.Range("F159").Value = 1

# San Bartolome / Sahara Occidental swap places in the ranking (row 215 <-> row 216)
$ws.Range("A215").Value = "Sahara Occidental"
$ws.Range("A216").Value = "San Bartolome"
